# Recalculate HK_R_acc_G accuracy values for the Global Thresholding sheet.
# The first row holds the series label; rows 2-50 hold the per-run accuracy
# values used to compute the mean accuracy for this thresholding method.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assert the header label (kept identical, just re-written so the
# workbook's shared-string table is refreshed along with the data below).
$ws.Cells.Item(1, 1).Value = "HK_R_acc_G"

$ws.Cells.Item(2, 1).Value = 47.668955547524391
$ws.Cells.Item(3, 1).Value = 47.668955547524391
$ws.Cells.Item(4, 1).Value = 46.114925912540656
$ws.Cells.Item(5, 1).Value = 44.091073364654861
$ws.Cells.Item(6, 1).Value = 44.452475605348752
$ws.Cells.Item(7, 1).Value = 45.247560534875312
$ws.Cells.Item(8, 1).Value = 47.524394651246837
$ws.Cells.Item(9, 1).Value = 48.138778460426458
$ws.Cells.Item(10, 1).Value = 47.596675099385614
$ws.Cells.Item(11, 1).Value = 47.632815323455006
$ws.Cells.Item(12, 1).Value = 51.680520419226596
$ws.Cells.Item(13, 1).Value = 51.535959522949035
$ws.Cells.Item(14, 1).Value = 47.45211420310806
$ws.Cells.Item(15, 1).Value = 48.174918684495843
$ws.Cells.Item(16, 1).Value = 47.668955547524391
$ws.Cells.Item(17, 1).Value = 46.946151066136608
$ws.Cells.Item(18, 1).Value = 48.102638236357066
$ws.Cells.Item(19, 1).Value = 47.994217564148897
$ws.Cells.Item(20, 1).Value = 48.391760028912181
$ws.Cells.Item(21, 1).Value = 46.620889049512101
$ws.Cells.Item(22, 1).Value = 48.50018070112035
$ws.Cells.Item(23, 1).Value = 49.222985182508133
$ws.Cells.Item(24, 1).Value = 53.234550054210338
$ws.Cells.Item(25, 1).Value = 52.909288037585831
$ws.Cells.Item(26, 1).Value = 47.560534875316229
$ws.Cells.Item(27, 1).Value = 47.560534875316229
$ws.Cells.Item(28, 1).Value = 47.488254427177452
$ws.Cells.Item(29, 1).Value = 53.089989157932784
$ws.Cells.Item(30, 1).Value = 52.620166245030717
$ws.Cells.Item(31, 1).Value = 47.415973979038675
$ws.Cells.Item(32, 1).Value = 49.222985182508133
$ws.Cells.Item(33, 1).Value = 45.681243223707988
$ws.Cells.Item(34, 1).Value = 45.861944344054933
$ws.Cells.Item(35, 1).Value = 48.355619804842789
$ws.Cells.Item(36, 1).Value = 48.644741597397903
$ws.Cells.Item(37, 1).Value = 52.728586917238886
$ws.Cells.Item(38, 1).Value = 49.042284062161187
$ws.Cells.Item(39, 1).Value = 48.717022045536687
$ws.Cells.Item(40, 1).Value = 49.620527647271409
$ws.Cells.Item(41, 1).Value = 48.102638236357066
$ws.Cells.Item(42, 1).Value = 48.24719913263462
$ws.Cells.Item(43, 1).Value = 48.24719913263462
$ws.Cells.Item(44, 1).Value = 47.162992410552945
$ws.Cells.Item(45, 1).Value = 47.632815323455006
$ws.Cells.Item(46, 1).Value = 46.620889049512101
$ws.Cells.Item(47, 1).Value = 47.090711962414169
$ws.Cells.Item(48, 1).Value = 45.93422479219371
$ws.Cells.Item(49, 1).Value = 47.162992410552945
